# docs/epexspot_prices.xlsx - "Gaz" sheet update
#
# The daily price feed was missing two days (2025-06-21 and 2025-06-22).
# They need to be inserted, in date order, right before the existing
# 2025-06-23 row - pushing the two rows that already follow it
# (2025-06-23 and 2025-06-25) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaz")

# 1) Capture the two rows that currently occupy rows 7-8
#    (2025-06-23 / 40.9 and 2025-06-25 / 34.75) before they get overwritten.
$date7 = $ws.Range("A7").Text
$price7 = $ws.Range("B7").Value2
$date8 = $ws.Range("A8").Text
$price8 = $ws.Range("B8").Value2

# 2) Move them down to rows 9-10 to make room for the two new rows.
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = $date7
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = $price7

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = $date8
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = $price8

# 3) Write the two newly-available daily prices into rows 7-8.
#    NumberFormat is forced to Text first so Excel doesn't reinterpret the
#    ISO date string as a date serial; Style is reset back to "Normal"
#    afterwards so no extra formatting is left behind on the cell.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-06-21"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 40.275

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025-06-22"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = 40.275
